$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column I: "RESPONSE_PHONE" header and its value for the data row.
# Leave these cells at their default (General) format, which matches the
# target formatting semantics for the new column.
$ws.Range("I1").Value = "RESPONSE_PHONE"
$ws.Range("I2").Value = "1-570-236-7033"

# Column H ("CODE") reverts to the same plain text ("@") format used by the
# other string columns (e.g. column G) instead of its previous dedicated
# text style.
$ws.Range("H1").NumberFormat = "@"
$ws.Range("H2").NumberFormat = "@"

# Give the new column I its own width (~17.37 characters).
$ws.Columns.Item(9).ColumnWidth = 16.45

# Move the active selection to match where editing left off.
$ws.Range("I5").Select()
